$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the cell values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the shared "boxed / bold / centered" style once on B1 ...
$b1 = $ws.Range("B1")
$b1.Font.Bold = $true
$b1.HorizontalAlignment = -4108  # xlCenter
$b1.VerticalAlignment = -4160    # xlTop
$b1.Borders.Weight = 2           # xlThin

# ... then clone the exact same format onto A2 via a format-only paste
# (re-running the same sequence of property mutations on a second cell
# would otherwise leave a stray, unused style behind).
$b1.Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
